# 141: 31/12 09:51 LP1912+6203+6173
# New scrape timestamp that this commit records.
$newTimestamp = "31/12/2025 06:51:36"

$wb = $excel.ActiveWorkbook

# --- Sheet "LP1912" (sheet1): 10 new trips appended (rows 694-703) ---
$wsLP = $wb.Worksheets.Item("LP1912")

# Row 2 col A: "Última actualización: ..." timestamp bump.
$wsLP.Cells.Item(2, 1).Value = "Última actualización: $newTimestamp"
# Row 3 col A: "Total filas: ..." bumped from 692 to 702 (+10 new rows).
$wsLP.Cells.Item(3, 1).Value = "Total filas: 702"

$lp1912NewRows = @(
    @("06:51:26", "06:54", "14_ABASTO",      3, "LP1912", "31/12/2025"),
    @("06:51:26", "07:01", "16_SANTA ANA",  10, "LP1912", "31/12/2025"),
    @("06:51:26", "07:16", "16_SANTA ANA",  25, "LP1912", "31/12/2025"),
    @("06:51:26", "07:29", "14_ABASTO",     38, "LP1912", "31/12/2025"),
    @("06:51:26", "07:37", "27_EL RETIRO",  46, "LP1912", "31/12/2025"),
    @("06:51:26", "07:51", "15_ABASTO",     60, "LP1912", "31/12/2025"),
    @("06:51:26", "08:03", "17_ROMERO",     72, "LP1912", "31/12/2025"),
    @("06:51:26", "08:03", "23_HERNANDEZ",  72, "LP1912", "31/12/2025"),
    @("06:51:26", "08:14", "10_OLMOS",      83, "LP1912", "31/12/2025"),
    @("06:51:26", "08:30", "14_ABASTO",     99, "LP1912", "31/12/2025")
)

$startRow = 694
for ($i = 0; $i -lt $lp1912NewRows.Count; $i++) {
    $r = $startRow + $i
    $row = $lp1912NewRows[$i]

    # Column A (Hora_Scrap header says otherwise, but col A on this sheet
    # only ever carries the "Última actualización"/"Total filas" banner text
    # on rows 2-3; for data rows it is blank) - write as explicit empty text
    # (not a truly-empty/never-touched cell) to match the source export.
    $wsLP.Cells.Item($r, 1).Value = "'"
    $wsLP.Cells.Item($r, 1).Style = "Normal"

    $wsLP.Cells.Item($r, 2).Value = $row[0]
    $wsLP.Cells.Item($r, 3).Value = $row[1]
    $wsLP.Cells.Item($r, 4).Value = $row[2]
    $wsLP.Cells.Item($r, 5).Value = $row[3]
    $wsLP.Cells.Item($r, 6).Value = $row[4]
    $wsLP.Cells.Item($r, 7).Value = $row[5]
}

# --- Sheet "LP1912-215" (sheet2): timestamp bump only, no new rows ---
$wsLP215 = $wb.Worksheets.Item("LP1912-215")
$wsLP215.Cells.Item(2, 1).Value = "Última actualización: $newTimestamp"

# --- Sheet "6203-6173" (sheet3): 2 new trips appended (rows 85-86) ---
$ws6203 = $wb.Worksheets.Item("6203-6173")
$ws6203.Cells.Item(2, 1).Value = "Última actualización: $newTimestamp"
$ws6203.Cells.Item(3, 1).Value = "Total filas: 85"

$sheet3NewRows = @(
    @("31/12/2025", "06:51:36", "07:27", "215A_LA PLATA", 36, "L6173"),
    @("31/12/2025", "06:51:36", "08:10", "215A_LA PLATA", 79, "L6173")
)

$startRow3 = 85
for ($i = 0; $i -lt $sheet3NewRows.Count; $i++) {
    $r = $startRow3 + $i
    $row = $sheet3NewRows[$i]

    $ws6203.Cells.Item($r, 1).Value = "'"
    $ws6203.Cells.Item($r, 1).Style = "Normal"

    $ws6203.Cells.Item($r, 2).Value = $row[0]
    $ws6203.Cells.Item($r, 3).Value = $row[1]
    $ws6203.Cells.Item($r, 4).Value = $row[2]
    $ws6203.Cells.Item($r, 5).Value = $row[3]
    $ws6203.Cells.Item($r, 6).Value = $row[4]
    $ws6203.Cells.Item($r, 7).Value = $row[5]
}
